# Remove the incomplete "INV-5678" entry (row 2) from the invoices sheet.
# This shifts the remaining rows up by one, and the "Processed At" timestamps
# for the shifted rows are refreshed to reflect the new processing pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 2 (INV-5678 / Incomplete), shifting rows 3-5 up to 2-4.
$ws.Rows.Item(2).Delete()

# Refresh the "Processed At" column (E) timestamps for the shifted rows.
$ws.Range("E2").Value = "2026-02-03 19:18:04"
$ws.Range("E3").Value = "2026-02-03 19:18:05"
$ws.Range("E4").Value = "2026-02-03 19:18:05"
